# Add "Aquila Optimizer" (AO) as entry #38 of the Swarm-based group.
# Before: row 46 is the last (37th) Swarm entry and is blank past column B
#         (B46 = 37), row 47 starts the Physics group.
# After:  row 46 gets the new AO data, a new row 47 (B47 = 38) closes out
#         the Swarm group's STT counter, and everything from the old row 47
#         onward shifts down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row right after the current last Swarm row (46); this pushes
# the old row 47 (start of the Physics group) down to row 48, etc.
$ws.Rows.Item(47).Insert()

# Populate the new Swarm entry on (now fully blank-past-B) row 46.
$ws.Range("C46").Value = "Aquila Optimizer"
$ws.Range("D46").Value = "AO"
$ws.Range("E46").Value = 2021
$ws.Range("F46").Value = "original"
$ws.Range("G46").Value = "no"
$ws.Range("H46").Value = "yes"
$ws.Range("I46").Value = "strong"
$ws.Range("J46").Value = "yes"
$ws.Range("K46").Value = 2
$ws.Range("L46").Value = "easy"

# The freshly-inserted row 47 just continues the Swarm group's STT counter.
$ws.Range("B47").Value = 38

# Restore the view: scroll so row 73 is the top row, and leave the
# selection on C47 (matches the author's on-save cursor position).
$excel.Goto($ws.Range("A73"), $true)
$ws.Range("C47").Select()
